$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# Insert a new row at position 83 (pushes e060 and everything after it down by one)
$ws.Rows("83:83").Insert()

# New row 83: e056 label + body text (added first so the new shared strings land
# right after the existing e055 label string, before the (soon to be updated) e055 body)
$ws.Range("A83").Value = "e056"

$e056Text = "<Bold>e056 Repair Gun</Bold> `n<InlineUIContainer><Button Content='r4.74.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   `n<LineBreak/><LineBreak/>`nAttempt to repair malfunction gun by rolling on the <InlineUIContainer><Button Content=Gun Malfunction' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table.`n<LineBreak/><LineBreak/>"
$ws.Range("B83").Value = $e056Text

# Update B82 (e055 body text): the phrase "out of PERISCOPE_REPLACEMENT" right before
# ". Click image to continue." is removed.
$e055Text = "<Bold>e055 Replace Periscopes</Bold> `n<InlineUIContainer><Button Content='r4.74.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   `n<LineBreak/><LineBreak/>`nReplacing PERISCOPE_REPLACEMENT out of PERISCOPE_REPLACEMENT_TOTAL left as shown on the After Action Report `n<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. Click image to continue.`n<LineBreak/><LineBreak/>`n                                            <InlineUIContainer><Image Name='BrokenPeriscope' Height='350' Width='222'></Image></InlineUIContainer>"
$ws.Range("B82").Value = $e055Text

$ws.Rows("83:83").RowHeight = 75

# Update selection / view state to match the target workbook
$ws.Application.ActiveWindow.ScrollRow = 80
$ws.Range("B82").Select()
